$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Gender" header in K1, copying the formatting from the other
# header cells (e.g. J1) so it picks up the same style.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Gender"

# Gender column formula: derive M/F from the first letter of the
# category code in column D. Applied in the same row groupings as the
# existing column J shared formulas (K2 alone, then K3:K66, then
# K67:K111) so the fill matches how the sheet was actually built up.
$ws.Range("K2").Formula = '=IF(LEFT(D2,1)="M","M","F")'
$ws.Range("K3:K66").Formula = '=IF(LEFT(D3,1)="M","M","F")'
$ws.Range("K67:K111").Formula = '=IF(LEFT(D67,1)="M","M","F")'

$excel.CutCopyMode = $false
